$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Feb" sheet (sheet11.xml): a new SA ("AMISHA SAWATKAR") was inserted into
# the alphabetically-sorted table, shifting every row below it down by one,
# and several counters were refreshed for other reps.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Feb")

# Header: "Division" -> "Divfision"
$ws.Range("F1").Value = "Divfision"

# Make room for the new row - shifts old rows 3..14 down to 4..15, carrying
# their existing cell content/types along automatically.
$ws.Rows.Item(3).Insert()

# Row 2 (ABHAYSINGH DESAI): Link Triggered count refreshed
$ws.Range("B2").Value = 7

# Row 3: brand-new rep inserted here
$ws.Range("A3").Value = "AMISHA SAWATKAR"
$ws.Range("B3").Value = 1
# C3/D3 stay blank (like the other reps with no Response/Concern Count data) -
# use the text quote-prefix trick so the cell is stored as an empty text
# value (matching the rest of the blank cells in this column) instead of
# being left out of the sheet entirely, then strip the quote-prefix style
# back off so no visible formatting change remains.
$ws.Range("C3").Formula = "'"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Formula = "'"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "NAGPUR_KAMPTHEE ROAD"

# Row 4 (ANIKET DESHMUKH, shifted from old row 3): refreshed count
$ws.Range("B4").Value = 44

# Row 5 (ANIKET FUSATE, shifted from old row 4): unchanged, no edits needed

# Row 6 (ASHUTOSH GURAV, shifted from old row 5): refreshed count
$ws.Range("B6").Value = 29

# Row 7 (AVINASH KAMBLE, shifted from old row 6): refreshed count + rate
$ws.Range("B7").Value = 38
$ws.Range("E7").Value = 52.63

# Row 8 (BHARAT GAVADE, shifted from old row 7): unchanged, no edits needed

# Row 9 (DNYANESHWAR GAWADE, shifted from old row 8): refreshed count
$ws.Range("B9").Value = 15

# Row 10 (JAVED RAMPURE, shifted from old row 9): refreshed count
$ws.Range("B10").Value = 13

# Row 11 (MANOJ PATIL, shifted from old row 10): refreshed count
$ws.Range("B11").Value = 16

# Row 12 (MOHSIN ALI, shifted from old row 11): refreshed count
$ws.Range("B12").Value = 21

# Row 13 (SANJAY RAMKELKAR, shifted from old row 12): refreshed count
$ws.Range("B13").Value = 9

# Row 14 (UJJWAL MAHAJAN, shifted from old row 13): refreshed count
$ws.Range("B14").Value = 37

# Row 15 (VAIBHAV PANCHAL, shifted from old row 14): refreshed count/rate/division
$ws.Range("B15").Value = 41
$ws.Range("E15").Value = 24.39
$ws.Range("F15").Value = "YAVATMAL"

# ---------------------------------------------------------------------------
# "Dec" sheet (sheet9.xml): just a selection change (F3)
# ---------------------------------------------------------------------------
$wsDec = $wb.Worksheets.Item("Dec")
$wsDec.Range("F3").Select()

# Re-activate "Feb" (it was, and remains, the selected tab) and move its
# selection from A7 to F7.
$ws.Activate()
$ws.Range("F7").Select()
